$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New measurement added for the 10000-iterations convolution run (rows 7 and 11)
$ws.Range("A7").Value = 10000
$ws.Range("A11").Value = 10000

# Updated timing results (multiprocessing numbers)
$ws.Range("F5").Value = 18.975380000000001
$ws.Range("F7").Value = 9.1827000000000005
$ws.Range("F9").Value = 21.533545
$ws.Range("F11").Value = 10.564019999999999

# Move the active selection to A11, matching where the author left the cursor
$ws.Range("A11").Select()
